$d = $word.ActiveDocument

# 1) Merge "12 - 13 неделя: тестирование " + "мобильной" + " версии сервиса" into one run
$d.Content.Find.Execute("12 " + [char]8212 + " 13 неделя: тестирование " + [char]2 + " версии сервиса", $false, $true, $false, $false, $false, $true, 1, $false, "12 " + [char]8212 + " 13 неделя: тестирование мобильной версии сервиса", 2) | Out-Null

# 2) Merge "13 неделя: исправление багов в " + "мобильной" + " версии сервиса"
$d.Content.Find.Execute("13 неделя: исправление багов в " + [char]2 + " версии сервиса", $false, $true, $false, $false, $false, $true, 1, $false, "13 неделя: исправление багов в мобильной версии сервиса", 2) | Out-Null

# 3) "6" + " неделя" -> "6 неделя" (careful, match only the heading occurrence)
$d.Content.Find.Execute([char]2 + " неделя", $false, $true, $false, $false, $false, $true, 1, $false, "6 неделя", 2) | Out-Null

# 4) "Составление тест-плана" -> "Составление тест-требований"
$d.Content.Find.Execute("Составление тест-плана", $false, $true, $false, $false, $false, $true, 1, $false, "Составление тест-требований", 2) | Out-Null

Write-Output "done"
